# Generate Report for Handback
#
# The localization status report is refreshed after the de-de handback
# completes: the "Ready for handoff" status becomes "Handed back: in sync
# with en-US" everywhere it is shown (Overview rollup + the zh-cn/de-de
# detail sheets), the "Latest Handback DateTime" stamps advance to the
# handback run's timestamps, and the stale "version mismatch" Error Detail
# text is cleared now that the handback is in sync.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Column widths in this report grow/shrink to fit the new Status / Error
# Detail text; ColumnWidth is specified in "characters" and gets quantized
# by the host to whole pixels, so these are the closest settable values to
# the report's rendered widths.
$statusColWidth = 29.166666666666668
$errorColWidth = 12.833333333333334

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn (E) / de-de (F) status columns ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E1").ColumnWidth = $statusColWidth
$wsOverview.Range("F1").ColumnWidth = $statusColWidth

# --- zh-cn detail sheet ---
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-08-12 09:00:05"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("C1").ColumnWidth = $statusColWidth
$wsZhCn.Range("P1").ColumnWidth = $errorColWidth

# --- de-de detail sheet ---
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-08-12 09:00:26"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("C1").ColumnWidth = $statusColWidth
$wsDeDe.Range("P1").ColumnWidth = $errorColWidth
